$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (best-effort; engine quantizes ColumnWidth to 1/6-character
#     buckets, so these are the closest reachable values to the target
#     4.5546875 / 7.6640625 / 8.77734375 / 12.21875 stored widths) ---
$ws.Columns.Item(1).ColumnWidth = 3.5833333333333335
$ws.Columns.Item(2).ColumnWidth = 6.75
$ws.Columns.Item(3).ColumnWidth = 7.916666666666667
$ws.Columns.Item(4).ColumnWidth = 11.25

# --- Rows 2-8: clear the trial data back to placeholder zero rows ---
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
    $ws.Cells.Item($r, 2).Value = ""
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}

# --- Rows 79-81: fill in the completed trial rows at the end of the log ---
$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = "T7-5.jpg"
$ws.Cells.Item(79, 3).Value = ""
$ws.Cells.Item(79, 4).Value = ""

$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "R7-9.jpg"
$ws.Cells.Item(80, 3).Value = ""
$ws.Cells.Item(80, 4).Value = ""

$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "R7-8.jpg"
$ws.Cells.Item(81, 3).Value = ""
$ws.Cells.Item(81, 4).Value = ""
